$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.468.05"
$ws.Range("E2").Value = '  -1.12%  '
$ws.Range("D3").Value = "'2.374.21"
$ws.Range("E3").Value = '  +4.75%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'235.35"
$ws.Range("E5").Value = '  +0.62%  '
$ws.Range("D6").Value = "'0.651"
$ws.Range("E6").Value = '  -0.71%  '
$ws.Range("D7").Value = "'72.54"
$ws.Range("E7").Value = '  +13.99%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = "'0.478"
$ws.Range("E9").Value = '  +5.58%  '
$ws.Range("D10").Value = "'0.0976"
$ws.Range("E10").Value = '  -0.20%  '
$ws.Range("D11").Value = "'56.91"
$ws.Range("E11").Value = '  -2.10%  '
$ws.Range("D12").Value = "'27.13"
$ws.Range("E12").Value = '  +1.72%  '
$ws.Range("D13").Value = "'2.732.61"
$ws.Range("E13").Value = '  +5.00%  '
$ws.Range("E14").Value = '  +0.41%  '
$ws.Range("E15").Value = '  +2.59%  '
$ws.Range("E16").Value = '  +2.50%  '
$ws.Range("D17").Value = "'0.857"
$ws.Range("E17").Value = '  +1.54%  '
$ws.Range("D18").Value = "'2.379.13"
$ws.Range("E18").Value = '  +5.18%  '
$ws.Range("D19").Value = "'43.414.41"
$ws.Range("E19").Value = '  -1.05%  '
$ws.Range("D20").Value = "'0.0₃0996"
$ws.Range("E20").Value = '  +1.17%  '
$ws.Range("D21").Value = "'6.35"
$ws.Range("E21").Value = '  +2.65%  '
$ws.Range("D22").Value = "'74.62"
$ws.Range("E22").Value = '  +0.97%  '
$ws.Range("D23").Value = "'250.81"
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("B25").Value = 'WEMIXToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D25").Value = "'3.69"
$ws.Range("E25").Value = '  +8.86%  '
$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").Value = "'2.47"
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("D27").Value = "'10.04"
$ws.Range("E27").Value = '  +1.07%  '
$ws.Range("E28").Value = '  -3.64%  '
$ws.Range("E29").Value = '  +2.32%  '
$ws.Range("D30").Value = "'174.21"
$ws.Range("E30").Value = '  +0.13%  '
$ws.Range("E31").Value = '  +6.11%  '
$ws.Range("E32").Value = '  -5.16%  '
$ws.Range("E33").Value = '  +0.13%  '
$ws.Range("E34").Value = '  +0.97%  '
$ws.Range("E35").Value = '  +0.72%  '
$ws.Range("D36").Value = "'5.07"
$ws.Range("E36").Value = '  +2.22%  '
$ws.Range("E37").Value = '  +7.15%  '
$ws.Range("E38").Value = '  +2.59%  '
$ws.Range("E39").Value = '  -0.95%  '
$ws.Range("E40").Value = '  +0.23%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = "'8.94"
$ws.Range("E41").Value = '  +2.37%  '
$ws.Range("B42").Value = 'BinanceUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("D43").Value = "'18.56"
$ws.Range("E43").Value = '  +7.37%  '
$ws.Range("D44").Value = "'1.19"
$ws.Range("E44").Value = '  +9.04%  '
$ws.Range("D45").Value = "'100.44"
$ws.Range("E45").Value = '  +1.56%  '
$ws.Range("D46").Value = "'4.50"
$ws.Range("E46").Value = '  +1.42%  '
$ws.Range("E47").Value = '  +1.67%  '
$ws.Range("D48").Value = "'0.0958"
$ws.Range("E48").Value = '  +0.65%  '
$ws.Range("D49").Value = "'1.451.90"
$ws.Range("E49").Value = '  -0.31%  '
$ws.Range("D50").Value = "'2.602.51"
$ws.Range("E50").Value = '  +5.08%  '
$ws.Range("E51").Value = '  -7.23%  '
